# Tweak spacing/size of the logo on the slide layout:
#  - "Picture 6" (the NWMSU logo image) is nudged right and shrunk a touch.
#  - "Content Placeholder 33" (the logo placeholder box behind it) gets the
#    matching new height so it still lines up with the picture.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$cl = $s.CustomLayout

$logo = $cl.Shapes.Item("Picture 6")
$logo.Left   = 3106.8
$logo.Top    = 40.32000165
$logo.Width  = 308.7159882019685
$logo.Height = 349.2

$ph = $cl.Shapes.Item("Content Placeholder 33")
$ph.Height = 349.2
